{"js": "// Apply the \"Added many more features\" edit to the Annihilator review doc.\n// Replace each old text with its new text using body.search() + insertText(replace),\n// which preserves the run formatting (bold/italic) of the matched text.\n\nconst body = context.document.body;\n\nconst replacements = [\n  {\n    oldText: \"Play Annihilator Slot for Free - Review of Annihilator Slot Machine\",\n    newText: \"Play Annihilator Slot Game Free\"\n  },\n  {\n    oldText: \"Expanding wilds and scatter symbols\",\n    newText: \"Exciting gameplay features like expanding wilds and bonus rounds\"\n  },\n  {\n    oldText: \"Lucrative bonus features\",\n    newText: \"Impressive maximum payout of 5,000 times your wager\"\n  },\n  {\n    oldText: \"Detailed graphics and atmosphere\",\n    newText: \"Meticulously detailed graphics and eerie atmosphere\"\n  },\n  {\n    oldText: \"Customizable soundtrack\",\n    newText: \"Customizable soundtrack featuring the band's greatest hits\"\n  },\n  {\n    oldText: \"Limited number of paylines\",\n    newText: \"Limited number of paylines (ten)\"\n  },\n  {\n    oldText: \"Higher minimum bet compared to other slots\",\n    newText: \"High wagering range might not be suitable for all players\"\n  },\n  {\n    oldText: \"Play Annihilator slot machine for free and read our review. Explore the game's features, graphics, soundtrack, theme, symbols, payouts, and more.\",\n    newText: \"Read our review of Annihilator, a slot game with exciting features and a free play option.\"\n  }\n];\n\nfor (const { oldText, newText } of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edit to the Annihilator review doc.\n# For each old/new text pair, use Find to locate the text and replace it by\n# assigning the matched Range's Text property directly (this avoids Word's\n# smart-quote autocorrection that Find.Execute's Replacement.Text triggers).\n\nfunction Replace-AllOccurrences {\n    param(\n        $Doc,\n        [string]$OldText,\n        [string]$NewText\n    )\n\n    $rng = $Doc.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $OldText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Forward = $true\n    $find.Wrap = 0\n\n    $replacedCount = 0\n    $found = $find.Execute()\n    while ($found) {\n        $rng.Text = $NewText\n        $replacedCount = $replacedCount + 1\n\n        # Continue searching after the text we just inserted.\n        $rng.Collapse(0)\n        $rng.End = $Doc.Content.End\n        $found = $find.Execute()\n    }\n\n    return $replacedCount\n}\n\n$d = $word.ActiveDocument\n\nReplace-AllOccurrences $d \"Play Annihilator Slot for Free - Review of Annihilator Slot Machine\" \"Play Annihilator Slot Game Free\" | Out-Null\nReplace-AllOccurrences $d \"Expanding wilds and scatter symbols\" \"Exciting gameplay features like expanding wilds and bonus rounds\" | Out-Null\nReplace-AllOccurrences $d \"Lucrative bonus features\" \"Impressive maximum payout of 5,000 times your wager\" | Out-Null\nReplace-AllOccurrences $d \"Detailed graphics and atmosphere\" \"Meticulously detailed graphics and eerie atmosphere\" | Out-Null\nReplace-AllOccurrences $d \"Customizable soundtrack\" \"Customizable soundtrack featuring the band's greatest hits\" | Out-Null\nReplace-AllOccurrences $d \"Limited number of paylines\" \"Limited number of paylines (ten)\" | Out-Null\nReplace-AllOccurrences $d \"Higher minimum bet compared to other slots\" \"High wagering range might not be suitable for all players\" | Out-Null\nReplace-AllOccurrences $d \"Play Annihilator slot machine for free and read our review. Explore the game's features, graphics, soundtrack, theme, symbols, payouts, and more.\" \"Read our review of Annihilator, a slot game with exciting features and a free play option.\" | Out-Null\n"}
